$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Coin price/volume updates scraped by the bot. Prices such as "212.11" or
# "19.30" must stay plain text (matching the original inline-string cells),
# so numeric-looking values are written with a leading "'" quote prefix —
# Excel strips the quote and stores the text as-is without treating it as a number.

$ws.Range('D2').Value = '26.238.97'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.587.54'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''212.11'
$ws.Range('E5').Value = '  +1.55%  '
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('D9').Value = '''0.0607'
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').Value = '''19.30'
$ws.Range('E10').Value = '  -1.28%  '
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').Value = '1.812.25'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').Value = '1.595.07'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').Value = '''0.519'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').Value = '''64.21'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Value = '26.242.83'
$ws.Range('D18').Value = '0.0₃0726'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').Value = '''7.37'
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('D20').Value = '''212.99'
$ws.Range('E20').Value = '  +2.43%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = '''4.26'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').Value = '''8.99'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('D25').Value = '''143.49'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '''7.00'
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('E28').Value = '  -0.72%  '
$ws.Range('D29').Value = '''15.19'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('E31').Value = '  +1.27%  '
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').Value = '1.338.06'
$ws.Range('E33').Value = '  +4.66%  '
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('E37').Value = '  -5.42%  '
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').Value = '''0.823'
$ws.Range('E39').Value = '  +1.75%  '
$ws.Range('D40').Value = '''5.77'
$ws.Range('E40').Value = '  +3.55%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = '''0.930'
$ws.Range('E42').Value = '  -16.19%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''0.768'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '''2.14'
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('D45').Value = '1.723.80'
$ws.Range('E45').Value = '  +0.82%  '
$ws.Range('D46').Value = '''61.23'
$ws.Range('E46').Value = '  -1.83%  '
$ws.Range('D47').Value = '''85.75'
$ws.Range('E47').Value = '  -3.38%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '''1.48'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '''0.0978'
$ws.Range('E49').Value = '  -2.60%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.0501'
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '''0.998'
$ws.Range('E51').Value = '  -0.21%  '
